# Applies the OOXML diff:
#  1. "Import and process time series data by running" is split into
#     5 runs: "Import" | "s" | " and process" | "es" | " time series data by running"
#     (net visible text becomes "Imports and processes time series data by running")
#  2/3/4. The "/one_day_ops/ - contains..." / "/ten_day_ops/ - contains..." /
#     "/demands/ - contains..." paragraphs each had their separator run
#     ("/ - " or "/demands/ - ") and the following "contains scripts that
#     create graphs and value boxes" run (previously two separate <w:r>
#     elements) merged into a single run. (A sibling paragraph,
#     "/situational_awareness/ - contains...", was already a single run and
#     is left untouched.)
#
# NOTE on this interpreter's PowerShell quirks (worked out empirically by
# probing with small throwaway scripts):
#   - Parameters bind *positionally only* - named args (-Foo bar) and
#     default parameter values are silently ignored/blank, and [switch]
#     params never read as $true. So every helper below takes plain
#     positional args, and every call site supplies all of them.
#   - `Outer $a (Inner x y)` (a *parenthesized* nested call as an argument)
#     can silently mis-bind when $a is a COM object, dropping arguments.
#     Assigning the inner call's result to a variable first (or using the
#     `$( ... )` subexpression form) avoids this, so every nested call
#     below is pre-evaluated into a variable before being passed on.

$d = $word.ActiveDocument

# ---- helpers -------------------------------------------------------------

# Locate the Range spanning $text within $scopeRange (e.g. a single
# paragraph's Range). Returns $null if not found.
function Find-TextIn($scopeRange, $text) {
    $scan = $d.Range($scopeRange.Start, $scopeRange.End)
    $found = $scan.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return $null
    }
    return $d.Range($scan.Start, $scan.End)
}

# Replace the contents of $range with the given sequence of <w:r> run XML,
# keeping the enclosing paragraph's own <w:p>/<w:pPr> (paraId, rsids, indent,
# etc.) untouched since only the inner range is targeted by InsertXML.
function Set-RunsXml($range, $runsXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# Build a single <w:r> (sz/szCs 22, matching this doc's body text runs)
# carrying $t; $preserve ($true/$false) controls xml:space="preserve".
function Run22($t, $preserve) {
    $sp = ''
    if ($preserve) {
        $sp = ' xml:space="preserve"'
    }
    return '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t' + $sp + '>' + $t + '</w:t></w:r>'
}

# ---- 1. "Import and process time series data by running" ---------------

$r1 = Find-TextIn $d.Content "Import and process time series data by running"
if ($r1 -ne $null) {
    $run1a = Run22 "Import" $false
    $run1b = Run22 "s" $false
    $run1c = Run22 " and process" $true
    $run1d = Run22 "es" $false
    $run1e = Run22 " time series data by running" $true
    $runs1 = $run1a + $run1b + $run1c + $run1d + $run1e
    Set-RunsXml $r1 $runs1
}

# ---- 2/3/4. merge the separator run + "contains scripts..." run ---------
# Target by paragraph text so the edit is resilient to the exact paragraph
# index; each paragraph is handled independently (scoped Find), so the
# already-merged "situational_awareness" paragraph is never touched.

$mergedSlash = Run22 "/ - contains scripts that create graphs and value boxes" $false
$mergedDemands = Run22 "/demands/ - contains scripts that create graphs and value boxes" $false

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $ptext = $p.Range.Text

    if ($ptext -like "/one_day_ops/*contains scripts that create graphs and value boxes*") {
        $hit = Find-TextIn $p.Range "/ - contains scripts that create graphs and value boxes"
        if ($hit -ne $null) {
            Set-RunsXml $hit $mergedSlash
        }
    } elseif ($ptext -like "/ten_day_ops/*contains scripts that create graphs and value boxes*") {
        $hit = Find-TextIn $p.Range "/ - contains scripts that create graphs and value boxes"
        if ($hit -ne $null) {
            Set-RunsXml $hit $mergedSlash
        }
    } elseif ($ptext -like "/demands/*contains scripts that create graphs and value boxes*") {
        $hit = Find-TextIn $p.Range "/demands/ - contains scripts that create graphs and value boxes"
        if ($hit -ne $null) {
            Set-RunsXml $hit $mergedDemands
        }
    }
}

Write-Host "Edits applied"
